$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 (pushes the existing rows 12..54 down to 13..55)
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new weekly record
$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C12").Value = "Ñuble"
$ws.Range("D12").Value = 44859
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 100112001
$ws.Range("G12").Value = "Berenjena"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 13000
$ws.Range("M12").Value = 12500
$ws.Range("N12").Value = "$/caja 60 unidades"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 208
$ws.Range("Q12").Value = 60
$ws.Range("R12").Value = "Hortaliza"
